# Atualizei dados bibi e add
# faturamento_diario.xlsx - update June (mes=6) daily revenue figures and
# add two missing days (26 and 27) to the June block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1) Correct existing June (Mes=6) totals ---------------------------------
$ws.Range("B4").Value2  = 11872.2    # Dia 4  -> total_venda
$ws.Range("B15").Value2 = 5774.35    # Dia 20 -> total_venda
$ws.Range("B17").Value2 = 70494.86   # Dia 23 -> total_venda
$ws.Range("B18").Value2 = 15341.74   # Dia 24 -> total_venda
$ws.Range("B19").Value2 = 17254.32   # Dia 25 -> total_venda

# --- 2) Insert two new rows for June days 26 and 27 --------------------------
# They belong right after the existing Dia=25 row (row 19), before the
# May block that currently starts at row 20.
$ws.Range("A20:E21").EntireRow.Insert()

$ws.Range("A20").Value2 = 26
$ws.Range("B20").Value2 = 21676.75
$ws.Range("C20").Value2 = 6
$ws.Range("D20").Value2 = 2025
$ws.Range("E20").Value2 = "06/2025"

$ws.Range("A21").Value2 = 27
$ws.Range("B21").Value2 = 504.45
$ws.Range("C21").Value2 = 6
$ws.Range("D21").Value2 = 2025
$ws.Range("E21").Value2 = "06/2025"
